$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for each data row (2..230).
# Update all of them from 45182 (2023-09-13) to 45184 (2023-09-15).
$ws.Range("C2:C230").Value2 = 45184
